# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Mayo de 2020 a las 01:35"

# Row 4 (Estados Unidos)
$ws.Range("B4").Value = 1455750
$ws.Range("C4").Value = 25403
$ws.Range("D4").Value = 316724
$ws.Range("E4").Value = 1052147
$ws.Range("F4").Value = 16240
$ws.Range("G4").Value = 1682
$ws.Range("H4").Value = 86879

# Row 17 (Canada)
$ws.Range("B17").Value = 73401
$ws.Range("C17").Value = 1123
$ws.Range("E17").Value = 31838

# Row 64 (Nigeria)
$ws.Range("B64").Value = 5162
$ws.Range("C64").Value = 191
$ws.Range("D64").Value = 1180
$ws.Range("E64").Value = 3815
$ws.Range("G64").Value = 3
$ws.Range("H64").Value = 167

# Row 97 (El Salvador)
$ws.Range("E97").Value = 684
$ws.Range("G97").Value = 3
$ws.Range("H97").Value = 23

# Row 139 (Vietnam)
$ws.Range("B139").Value = 312
$ws.Range("C139").Value = 24
$ws.Range("E139").Value = 52
